$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 253
$ws.Range("B253").Value = 6775597
$ws.Range("E253").Value = 45402.41666666666
$ws.Range("F253").Value = "Pogon Szczecin"
$ws.Range("G253").Value = "Piast Gliwice"
$ws.Range("K253").Value = 1.85
$ws.Range("L253").Value = 3.5
$ws.Range("M253").Value = 4.2
$ws.Range("N253").Value = 1.85
$ws.Range("O253").Value = 3.5
$ws.Range("P253").Value = 4.2
$ws.Range("Q253").Value = -0.5
$ws.Range("R253").Value = 1.85
$ws.Range("S253").Value = 2
$ws.Range("T253").Value = 2.25
$ws.Range("U253").Value = 2.025
$ws.Range("V253").Value = 1.825

# Row 254
$ws.Range("B254").Value = 6774472
$ws.Range("E254").Value = 45402.52083333334
$ws.Range("F254").Value = "Ruch Chorzow"
$ws.Range("G254").Value = "Widzew Lodz"
$ws.Range("K254").Value = 2.55
$ws.Range("L254").Value = 3.3
$ws.Range("M254").Value = 2.7
$ws.Range("N254").Value = 2.1
$ws.Range("O254").Value = 3.4
$ws.Range("P254").Value = 3.3
$ws.Range("Q254").Value = -0.25
$ws.Range("R254").Value = 1.875
$ws.Range("S254").Value = 1.975
$ws.Range("T254").Value = 2.5
$ws.Range("U254").Value = 1.975
$ws.Range("V254").Value = 1.875

# Row 255
$ws.Range("B255").Value = 6775594
$ws.Range("E255").Value = 45402.625
$ws.Range("F255").Value = "Zaglebie Lubin"
$ws.Range("G255").Value = "Jagiellonia Bialystok"
$ws.Range("K255").Value = 2.875
$ws.Range("L255").Value = 3.4
$ws.Range("M255").Value = 2.375
$ws.Range("N255").Value = 2.9
$ws.Range("O255").Value = 3.4
$ws.Range("P255").Value = 2.375
$ws.Range("Q255").Value = 0.25
$ws.Range("R255").Value = 1.8
$ws.Range("S255").Value = 2.05
$ws.Range("T255").Value = 2.75
$ws.Range("U255").Value = 2.025
$ws.Range("V255").Value = 1.825

# Row 256
$ws.Range("B256").Value = 6850054
$ws.Range("E256").Value = 45403.3125
$ws.Range("F256").Value = "Cracovia Krakow"
$ws.Range("G256").Value = "Puszcza Niepolomice"
$ws.Range("K256").Value = 1.7
$ws.Range("L256").Value = 3.8
$ws.Range("M256").Value = 4.75
$ws.Range("N256").Value = 1.7
$ws.Range("O256").Value = 3.8
$ws.Range("P256").Value = 4.75
$ws.Range("Q256").Value = -0.75
$ws.Range("R256").Value = 1.9
$ws.Range("S256").Value = 1.95
$ws.Range("T256").Value = 2.5
$ws.Range("U256").Value = 1.975
$ws.Range("V256").Value = 1.875

# Row 257
$ws.Range("B257").Value = 6830603
$ws.Range("E257").Value = 45403.41666666666
$ws.Range("F257").Value = "LKS Lodz"
$ws.Range("G257").Value = "Lech Poznan"
$ws.Range("K257").Value = 4.75
$ws.Range("L257").Value = 4
$ws.Range("M257").Value = 1.65
$ws.Range("N257").Value = 5.25
$ws.Range("O257").Value = 4
$ws.Range("P257").Value = 1.6
$ws.Range("Q257").Value = 1
$ws.Range("R257").Value = 1.825
$ws.Range("S257").Value = 2.025
$ws.Range("T257").Value = 2.75
$ws.Range("U257").Value = 1.975
$ws.Range("V257").Value = 1.875

# Row 258
$ws.Range("B258").Value = 6775596
$ws.Range("E258").Value = 45403.52083333334
$ws.Range("F258").Value = "Legia Warsaw"
$ws.Range("G258").Value = "Slask Wroclaw"
$ws.Range("K258").Value = 1.7
$ws.Range("L258").Value = 3.6
$ws.Range("M258").Value = 5
$ws.Range("N258").Value = 1.7
$ws.Range("O258").Value = 3.6
$ws.Range("P258").Value = 5
$ws.Range("Q258").Value = -0.75
$ws.Range("R258").Value = 1.925
$ws.Range("S258").Value = 1.925
$ws.Range("T258").Value = 2.5
$ws.Range("U258").Value = 2.05
$ws.Range("V258").Value = 1.8

# Row 259
$ws.Range("B259").Value = 6885526
$ws.Range("E259").Value = 45404.58333333334
$ws.Range("F259").Value = "Warta Poznan"
$ws.Range("G259").Value = "Stal Mielec"
$ws.Range("K259").Value = 2.1
$ws.Range("L259").Value = 3.1
$ws.Range("M259").Value = 3.8
$ws.Range("N259").Value = 2.1
$ws.Range("O259").Value = 3.1
$ws.Range("P259").Value = 3.8
$ws.Range("Q259").Value = -0.25
$ws.Range("R259").Value = 1.8
$ws.Range("S259").Value = 2.05
$ws.Range("T259").Value = 2
$ws.Range("U259").Value = 2.025
$ws.Range("V259").Value = 1.825

# Remove the two trailing rows that no longer have data (260 and 261)
$ws.Range("A260:A261").EntireRow.Delete()
